$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is the b.md entry ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-23 00:35:26"

# ---- zh-cn sheet: row 3 is the b.md entry ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-23 00:35:22"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d47d29851767a6fde02933654d918a0a7d1e468/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c91e32041ae1cf665728a714ec703ed5c2adce9/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 is the b.md entry ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-23 00:35:26"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d47d29851767a6fde02933654d918a0a7d1e468/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c91e32041ae1cf665728a714ec703ed5c2adce9/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17

Write-Output "Done"
